# Apply review-table update to Sheet1:
#  - shift/replace the review rows (2-15) with the refreshed review data
#  - add a new row (15) whose email/recovery columns use a wrap-text style
#  - add seven blank formatted rows (16-22) below the table
#  - move the active selection back to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the existing data rows (2-14) with the new review contents.
#    These rows keep their original per-column formatting, so only the
#    cell values need to change.
# ---------------------------------------------------------------------
$rows = @(
  @(2,  "com.hamxa.shaynachim", "bitcoin", "emmakrigel63@gmail.com ",  "sofershani9@gmail.com",    "27/5/2019 15:59", "excellent",                                               "no"),
  @(3,  "com.hamxa.shaynachim", "bitcoin", "erlichyotem@gmail.com ",   "rozend80@gmail.com",       "27/5/2019 15:59", "Things happen fast. Good app",                           "no"),
  @(4,  "com.hamxa.shaynachim", "bitcoin", "amramg25@gmail.com ",      "erlichyotem@gmail.com ",   "27/5/2019 15:59", "Bitcoin is a fast and furious technology. Great guide",  "no"),
  @(5,  "com.hamxa.shaynachim", "bitcoin", "fogelelad8@gmail.com",     "oamit1038@gmail.com",      "27/5/2019 15:59", "marvel app",                                              "no"),
  @(6,  "com.hamxa.shaynachim", "bitcoin", "mesikam455@gmail.com",     "imesika53@gmail.com",      "27/5/2019 15:59", "incredible",                                              "no"),
  @(7,  "com.hamxa.shaynachim", "bitcoin", "edenn0836@gmail.com",      "mesikam455@gmail.com",     "27/5/2019 15:59", "The way I wanted to read it",                            "no"),
  @(8,  "com.hamxa.shaynachim", "bitcoin", "frimanoren6@gmail.com",    "edenn0836@gmail.com",      "27/5/2019 15:59", "words of wisdom",                                         "no"),
  @(9,  "com.hamxa.shaynachim", "bitcoin", "goldfinshmulik@gmail.com", "frimanoren6@gmail.com",    "27/5/2019 15:59", "clear and easy",                                          "no"),
  @(10, "com.hamxa.shaynachim", "bitcoin", "elad86643@gmail.com",      "goldfinshmulik@gmail.com", "27/5/2019 15:59", "great app",                                               "no"),
  @(11, "com.hamxa.shaynachim", "bitcoin", "erezadmoni26@gmail.com",   "goldfinshmulik@gmail.com", "27/5/2019 15:59", "wise and simple",                                         "no"),
  @(12, "com.hamxa.shaynachim", "bitcoin", "margalitgal31@gmail.com",  "erezadmoni26@gmail.com",   "27/5/2019 15:59", "I cant believe it",                                       "yes"),
  @(13, "com.hamxa.shaynachim", "bitcoin", "maudanaor@gmail.com",      "margalitgal31@gmail.com",  "27/5/2019 15:59", "great great great",                                       "yes"),
  @(14, "com.hamxa.shaynachim", "bitcoin", "dmichal229@gmail.com",     "maudanaor@gmail.com",      "27/5/2019 15:59", "this app deserves 5 star",                                "yes")
)

foreach ($row in $rows) {
  $r = $row[0]
  $ws.Range("A$r").Value = $row[1]
  $ws.Range("B$r").Value = $row[2]
  $ws.Range("C$r").Value = $row[3]
  $ws.Range("D$r").Value = $row[4]
  $ws.Range("E$r").Value = $row[5]
  $ws.Range("F$r").Value = $row[6]
  $ws.Range("G$r").Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Add the brand-new row 15. Columns A and B reuse the formatting of
#    the other data rows; columns C and D get a new word-wrapping style
#    (same font as before, horizontal alignment back to general).
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").HorizontalAlignment = 1
$ws.Range("C15").WrapText = $true
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A15").Value = "com.hamxa.shaynachim"
$ws.Range("B15").Value = "bitcoin"
$ws.Range("C15").Value = "efiamid9@gmail.com"
$ws.Range("D15").Value = "maudanaor@gmail.com"
$ws.Range("E15").Value = "27/5/2019 15:59"
$ws.Range("F15").Value = "use it as it is"
$ws.Range("G15").Value = "yes"

$ws.Rows.Item(15).RowHeight = 13.8

# ---------------------------------------------------------------------
# 3) Add seven trailing, otherwise-empty rows (16-22) that only carry
#    the email/recovery column formatting forward.
# ---------------------------------------------------------------------
$ws.Range("C2").Copy()
for ($r = 16; $r -le 22; $r++) {
  $ws.Range("C$r").PasteSpecial(-4122)
  $ws.Range("D$r").PasteSpecial(-4122)
  $ws.Rows.Item($r).RowHeight = 13.8
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Restore the selection to A2 (top-left of the data, scrolled home).
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
$ws.Range("A2").Select()
